# Update localization status text: "Ready for handoff" -> "In Translation"
# across the Overview, zh-cn and de-de sheets, then shrink the affected
# status columns to fit the new (shorter) text, same as Excel's
# AutoFit/column-width update would do after a report refresh.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: zh-cn / de-de status columns are E and F (rows 2-3)
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# zh-cn sheet: Status column is C (rows 2-3)
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

# de-de sheet: Status column is C (rows 2-3)
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# Resize the status columns so they fit the new, shorter text (the
# report's generated column widths narrow along with the content).
$wsOverview.Columns("E:F").ColumnWidth = 12.5
$wsZhCn.Columns("C:C").ColumnWidth = 12.5
$wsDeDe.Columns("C:C").ColumnWidth = 12.5

Write-Output "Updated status text and resized columns"
